$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 25824.285  # H3: was 27222.8
$ws.Cells.Item(3, 10).Value = 25824.285  # J3: was 27222.8
$ws.Cells.Item(3, 12).Value = 25824.285  # L3: was 27222.8
$ws.Cells.Item(3, 14).Value = -26052.285  # N3: was -27450.8
$ws.Cells.Item(39, 8).Value = 992468.0600000001  # H39: was 882208.3
$ws.Cells.Item(39, 9).Value = 1323036.6  # I39: was 1221261.5
$ws.Cells.Item(39, 10).Value = 762.5  # J39: was 670
$ws.Cells.Item(39, 11).Value = 3969109.8  # K39: was 3663784.5
$ws.Cells.Item(39, 12).Value = 2287.5  # L39: was 2010
$ws.Cells.Item(39, 13).Value = -3968813.8  # M39: was -3663488.5
$ws.Cells.Item(39, 14).Value = -2879.5  # N39: was -2602
$ws.Cells.Item(64, 8).Value = 1224867.2  # H64: was 789063.8
$ws.Cells.Item(64, 9).Value = 2751500.5  # I64: was 3667634
$ws.Cells.Item(64, 10).Value = 3560.6  # J64: was 3999.182
$ws.Cells.Item(64, 11).Value = 2751500.5  # K64: was 3667634
$ws.Cells.Item(64, 12).Value = 3560.6  # L64: was 3999.182
$ws.Cells.Item(64, 13).Value = -2751252.5  # M64: was -3667386
$ws.Cells.Item(64, 14).Value = -4056.6  # N64: was -4495.182
$ws.Cells.Item(67, 8).Value = 1224867.2  # H67: was 789063.8
$ws.Cells.Item(67, 9).Value = 2751500.5  # I67: was 3667634
$ws.Cells.Item(67, 10).Value = 3560.6  # J67: was 3999.182
$ws.Cells.Item(67, 11).Value = 2751500.5  # K67: was 3667634
$ws.Cells.Item(67, 12).Value = 3560.6  # L67: was 3999.182
$ws.Cells.Item(67, 13).Value = -2750642.5  # M67: was -3666776
$ws.Cells.Item(67, 14).Value = -5276.6  # N67: was -5715.182
$ws.Cells.Item(102, 8).Value = 25824.285  # H102: was 27222.8
$ws.Cells.Item(102, 10).Value = 25824.285  # J102: was 27222.8
$ws.Cells.Item(102, 12).Value = 25824.285  # L102: was 27222.8
$ws.Cells.Item(102, 14).Value = -32314.285  # N102: was -33712.8
$ws.Cells.Item(127, 8).Value = 2324.9456  # H127: was 2205.9824
$ws.Cells.Item(127, 9).Value = 483.4  # I127: was 550
$ws.Cells.Item(127, 10).Value = 2509.1  # J127: was 2266.2
$ws.Cells.Item(127, 11).Value = 1450.2  # K127: was 1650
$ws.Cells.Item(127, 12).Value = 7527.299999999999  # L127: was 6798.599999999999
$ws.Cells.Item(127, 13).Value = 3509.8  # M127: was 3310
$ws.Cells.Item(127, 14).Value = -17447.3  # N127: was -16718.6
$ws.Cells.Item(129, 8).Value = 4077.1562  # H129: was 4187.387
$ws.Cells.Item(129, 9).Value = 20619.4  # I129: was 17275.334
$ws.Cells.Item(129, 10).Value = 1013.7778  # J129: was 1046.28
$ws.Cells.Item(129, 11).Value = 61858.2  # K129: was 51826.00199999999
$ws.Cells.Item(129, 12).Value = 3041.3334  # L129: was 3138.84
$ws.Cells.Item(129, 13).Value = -56858.2  # M129: was -46826.00199999999
$ws.Cells.Item(129, 14).Value = -13041.3334  # N129: was -13138.84
$ws.Cells.Item(132, 8).Value = 3791428.5  # H132: was 3909928.5
$ws.Cells.Item(132, 9).Value = 4035760  # I132: was 4170283.5
$ws.Cells.Item(132, 10).Value = 4291.25  # J132: was 4602.5
$ws.Cells.Item(132, 11).Value = 12107280  # K132: was 12510850.5
$ws.Cells.Item(132, 12).Value = 12873.75  # L132: was 13807.5
$ws.Cells.Item(132, 13).Value = -12104750  # M132: was -12508320.5
$ws.Cells.Item(132, 14).Value = -17933.75  # N132: was -18867.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 2388.1765  # H63: was 2341.1177
$ws.Cells.Item(63, 9).Value = 1670  # I63: was 1645.3636
$ws.Cells.Item(63, 10).Value = 3414.1428  # J63: was 3616.6667
$ws.Cells.Item(63, 11).Value = 1670  # K63: was 1645.3636
$ws.Cells.Item(63, 12).Value = 3414.1428  # L63: was 3616.6667
$ws.Cells.Item(63, 13).Value = -984  # M63: was -959.3635999999999
$ws.Cells.Item(63, 14).Value = -4786.1428  # N63: was -4988.6667
$ws.Cells.Item(66, 8).Value = 2388.1765  # H66: was 2341.1177
$ws.Cells.Item(66, 9).Value = 1670  # I66: was 1645.3636
$ws.Cells.Item(66, 10).Value = 3414.1428  # J66: was 3616.6667
$ws.Cells.Item(66, 11).Value = 8350  # K66: was 8226.817999999999
$ws.Cells.Item(66, 12).Value = 17070.714  # L66: was 18083.3335
$ws.Cells.Item(66, 13).Value = -4918  # M66: was -4794.817999999999
$ws.Cells.Item(66, 14).Value = -23934.714  # N66: was -24947.3335
$ws.Cells.Item(86, 8).Value = 36863  # H86: was 30133
$ws.Cells.Item(86, 10).Value = 50152  # J86: was 49981
$ws.Cells.Item(86, 12).Value = 50152  # L86: was 49981
$ws.Cells.Item(86, 14).Value = -52524  # N86: was -52353
$ws.Cells.Item(89, 8).Value = 36863  # H89: was 30133
$ws.Cells.Item(89, 10).Value = 50152  # J89: was 49981
$ws.Cells.Item(89, 12).Value = 150456  # L89: was 149943
$ws.Cells.Item(89, 14).Value = -162312  # N89: was -161799

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 117400  # H20: was 117500
$ws.Cells.Item(20, 9).Value = 127075  # I20: was 145071.42
$ws.Cells.Item(20, 10).Value = 40000  # J20: was 21000
$ws.Cells.Item(20, 11).Value = 127075  # K20: was 145071.42
$ws.Cells.Item(20, 12).Value = 40000  # L20: was 21000
$ws.Cells.Item(20, 13).Value = -126828  # M20: was -144824.42
$ws.Cells.Item(20, 14).Value = -40494  # N20: was -21494

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 71156.375  # H31: was 73149.55499999999
$ws.Cells.Item(31, 9).Value = 63318.25  # I31: was 72263.57000000001
$ws.Cells.Item(31, 10).Value = 77128.28999999999  # J31: was 73713.37
$ws.Cells.Item(31, 11).Value = 63318.25  # K31: was 72263.57000000001
$ws.Cells.Item(31, 12).Value = 77128.28999999999  # L31: was 73713.37
$ws.Cells.Item(31, 13).Value = -63023.25  # M31: was -71968.57000000001
$ws.Cells.Item(31, 14).Value = -77718.28999999999  # N31: was -74303.37
$ws.Cells.Item(34, 8).Value = 71156.375  # H34: was 73149.55499999999
$ws.Cells.Item(34, 9).Value = 63318.25  # I34: was 72263.57000000001
$ws.Cells.Item(34, 10).Value = 77128.28999999999  # J34: was 73713.37
$ws.Cells.Item(34, 11).Value = 63318.25  # K34: was 72263.57000000001
$ws.Cells.Item(34, 12).Value = 77128.28999999999  # L34: was 73713.37
$ws.Cells.Item(34, 13).Value = -63116.25  # M34: was -72061.57000000001
$ws.Cells.Item(34, 14).Value = -77532.28999999999  # N34: was -74117.37
$ws.Cells.Item(62, 8).Value = 2611.75  # H62: was 2699.8
$ws.Cells.Item(62, 9).Value = 2000  # I62: was 0
$ws.Cells.Item(62, 10).Value = 2699.1428  # J62: was 2699.8
$ws.Cells.Item(62, 11).Value = 2000  # K62: was 0
$ws.Cells.Item(62, 12).Value = 2699.1428  # L62: was 2699.8
$ws.Cells.Item(62, 13).Value = -1376  # M62: was None
$ws.Cells.Item(62, 14).Value = -3947.1428  # N62: was -3947.8
$ws.Cells.Item(64, 8).Value = 0  # H64: was 43208.43
$ws.Cells.Item(64, 10).Value = 0  # J64: was 43208.43
$ws.Cells.Item(64, 12).Value = 0  # L64: was 43208.43
$ws.Cells.Item(64, 14).ClearContents()  # N64: remove (was -43704.43)
$ws.Cells.Item(65, 8).Value = 2611.75  # H65: was 2699.8
$ws.Cells.Item(65, 9).Value = 2000  # I65: was 0
$ws.Cells.Item(65, 10).Value = 2699.1428  # J65: was 2699.8
$ws.Cells.Item(65, 11).Value = 10000  # K65: was 0
$ws.Cells.Item(65, 12).Value = 13495.714  # L65: was 13499
$ws.Cells.Item(65, 13).Value = -6880  # M65: was None
$ws.Cells.Item(65, 14).Value = -19735.714  # N65: was -19739
$ws.Cells.Item(67, 8).Value = 0  # H67: was 43208.43
$ws.Cells.Item(67, 10).Value = 0  # J67: was 43208.43
$ws.Cells.Item(67, 12).Value = 0  # L67: was 43208.43
$ws.Cells.Item(67, 14).ClearContents()  # N67: remove (was -44924.43)
$ws.Cells.Item(106, 8).Value = 27486.625  # H106: was 24388.875
$ws.Cells.Item(106, 10).Value = 27486.625  # J106: was 24388.875
$ws.Cells.Item(106, 12).Value = 27486.625  # L106: was 24388.875
$ws.Cells.Item(106, 14).Value = -30010.625  # N106: was -26912.875
$ws.Cells.Item(132, 8).Value = 2353.1082  # H132: was 2665.5625
$ws.Cells.Item(132, 9).Value = 2275.653  # I132: was 2608.5
$ws.Cells.Item(132, 10).Value = 2504.92  # J132: was 2774.5
$ws.Cells.Item(132, 11).Value = 6826.958999999999  # K132: was 7825.5
$ws.Cells.Item(132, 12).Value = 7514.76  # L132: was 8323.5
$ws.Cells.Item(132, 13).Value = -4296.958999999999  # M132: was -5295.5
$ws.Cells.Item(132, 14).Value = -12574.76  # N132: was -13383.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 4677.9  # H3: was 5720
$ws.Cells.Item(3, 9).Value = 3539.8572  # I3: was 4854.2856
$ws.Cells.Item(3, 10).Value = 7333.3335  # J3: was 8750
$ws.Cells.Item(3, 11).Value = 10619.5716  # K3: was 14562.8568
$ws.Cells.Item(3, 12).Value = 22000.0005  # L3: was 26250
$ws.Cells.Item(3, 13).Value = -10507.5716  # M3: was -14450.8568
$ws.Cells.Item(3, 14).Value = -22224.0005  # N3: was -26474
$ws.Cells.Item(4, 8).Value = 56279000  # H4: was 83334070
$ws.Cells.Item(4, 9).Value = 2166667.8  # I4: was 121.42857
$ws.Cells.Item(4, 10).Value = 83335170  # J4: was 200001600
$ws.Cells.Item(4, 11).Value = 6500003.399999999  # K4: was 364.28571
$ws.Cells.Item(4, 12).Value = 250005510  # L4: was 600004800
$ws.Cells.Item(4, 13).Value = -6499891.399999999  # M4: was -252.28571
$ws.Cells.Item(4, 14).Value = -250005734  # N4: was -600005024
$ws.Cells.Item(95, 8).Value = 6000  # H95: was 0
$ws.Cells.Item(95, 10).Value = 6000  # J95: was 0
$ws.Cells.Item(95, 12).Value = 18000  # L95: was 0
$ws.Cells.Item(95, 14).Value = -22118  # N95: was None
$ws.Cells.Item(106, 8).Value = 2499.8572  # H106: was 2500
$ws.Cells.Item(106, 10).Value = 2499.8572  # J106: was 2500
$ws.Cells.Item(106, 12).Value = 7499.571599999999  # L106: was 7500
$ws.Cells.Item(106, 14).Value = -9391.571599999999  # N106: was -9392
$ws.Cells.Item(131, 8).Value = 856.42  # H131: was 867.61
$ws.Cells.Item(131, 9).Value = 0  # I131: was 730
$ws.Cells.Item(131, 10).Value = 856.42  # J131: was 869
$ws.Cells.Item(131, 11).Value = 0  # K131: was 2190
$ws.Cells.Item(131, 12).Value = 2569.26  # L131: was 2607
$ws.Cells.Item(131, 13).ClearContents()  # M131: remove (was 2850)
$ws.Cells.Item(131, 14).Value = -12649.26  # N131: was -12687
$ws.Cells.Item(141, 8).Value = 2854.889  # H141: was 3933
$ws.Cells.Item(141, 9).Value = 2349  # I141: was 2852.8
$ws.Cells.Item(141, 10).Value = 3866.6667  # J141: was 5733.3335
$ws.Cells.Item(141, 11).Value = 7047  # K141: was 8558.400000000001
$ws.Cells.Item(141, 12).Value = 11600.0001  # L141: was 17200.0005
$ws.Cells.Item(141, 13).Value = -1867  # M141: was -3378.400000000001
$ws.Cells.Item(141, 14).Value = -21960.0001  # N141: was -27560.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 105780.45  # H70: was 100981.336
$ws.Cells.Item(70, 10).Value = 8088.4287  # J70: was 7702.25
$ws.Cells.Item(70, 12).Value = 8088.4287  # L70: was 7702.25
$ws.Cells.Item(70, 14).Value = -8628.4287  # N70: was -8242.25
$ws.Cells.Item(73, 8).Value = 105780.45  # H73: was 100981.336
$ws.Cells.Item(73, 10).Value = 8088.4287  # J73: was 7702.25
$ws.Cells.Item(73, 12).Value = 8088.4287  # L73: was 7702.25
$ws.Cells.Item(73, 14).Value = -9960.4287  # N73: was -9574.25
$ws.Cells.Item(105, 8).Value = 43326.668  # H105: was 45952.5
$ws.Cells.Item(105, 10).Value = 43326.668  # J105: was 45952.5
$ws.Cells.Item(105, 12).Value = 43326.668  # L105: was 45952.5
$ws.Cells.Item(105, 14).Value = -50314.668  # N105: was -52940.5
$ws.Cells.Item(134, 8).Value = 27032.715  # H134: was 27048.857
$ws.Cells.Item(134, 10).Value = 27032.715  # J134: was 27048.857
$ws.Cells.Item(134, 12).Value = 81098.145  # L134: was 81146.571
$ws.Cells.Item(134, 14).Value = -86168.145  # N134: was -86216.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2592.5117  # H132: was 2748.75
$ws.Cells.Item(132, 9).Value = 3125.76  # I132: was 3164.64
$ws.Cells.Item(132, 10).Value = 1851.8889  # J132: was 2055.6
$ws.Cells.Item(132, 11).Value = 9377.280000000001  # K132: was 9493.92
$ws.Cells.Item(132, 12).Value = 5555.6667  # L132: was 6166.799999999999
$ws.Cells.Item(132, 13).Value = -6847.280000000001  # M132: was -6963.92
$ws.Cells.Item(132, 14).Value = -10615.6667  # N132: was -11226.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 803.1177  # H113: was 839
$ws.Cells.Item(113, 9).Value = 570  # I113: was 610
$ws.Cells.Item(113, 10).Value = 1136.1428  # J113: was 1125.25
$ws.Cells.Item(113, 11).Value = 1710  # K113: was 1830
$ws.Cells.Item(113, 12).Value = 3408.4284  # L113: was 3375.75
$ws.Cells.Item(113, 13).Value = 460  # M113: was 340
$ws.Cells.Item(113, 14).Value = -7748.428400000001  # N113: was -7715.75
$ws.Cells.Item(132, 8).Value = 1785.2533  # H132: was 1961.3823
$ws.Cells.Item(132, 9).Value = 1398.7258  # I132: was 1589.8334
$ws.Cells.Item(132, 10).Value = 3628.6924  # J132: was 3394.5
$ws.Cells.Item(132, 11).Value = 4196.1774  # K132: was 4769.5002
$ws.Cells.Item(132, 12).Value = 10886.0772  # L132: was 10183.5
$ws.Cells.Item(132, 13).Value = -1666.1774  # M132: was -2239.5002
$ws.Cells.Item(132, 14).Value = -15946.0772  # N132: was -15243.5
